# Updated cryptos list on Sun Dec 24 12:54:33 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns with the latest scrape.
#
# Several Price values are numeric-looking strings that must stay TEXT
# (e.g. "1.00", "0.0950", "21.50") so trailing/leading zeros survive - a
# plain .Value assignment would let Excel coerce them to numbers and lose
# the formatting. We force those via a leading quote-prefix (like typing
# '1.00 into the cell) and then reset .Style back to "Normal" so the only
# thing that changes is the cell's stored value, not its formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.775.01'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '2.291.67'
$ws.Range("E3").Value = '  -0.17%  '
$c = $ws.Range("D4")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '
$c = $ws.Range("D5")
$c.Value = "'112.93"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +15.83%  '
$c = $ws.Range("D6")
$c.Value = "'269.04"
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$c = $ws.Range("D7")
$c.Value = "'0.625"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E8").Value = '  +0.06%  '
$c = $ws.Range("D9")
$c.Value = "'0.619"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.49%  '
$c = $ws.Range("D10")
$c.Value = "'48.23"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +5.61%  '
$c = $ws.Range("D11")
$c.Value = "'0.0950"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.65%  '
$c = $ws.Range("D12")
$c.Value = "'9.06"
$c.Style = "Normal"
$ws.Range("E12").Value = '  +14.35%  '
$ws.Range("E13").Value = '  +0.04%  '
$c = $ws.Range("D14")
$c.Value = "'15.81"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = '2.633.21'
$ws.Range("E15").Value = '  -0.13%  '
$c = $ws.Range("D16")
$c.Value = "'0.849"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '2.283.17'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").Value = '43.663.78'
$c = $ws.Range("D19")
$c.Value = "'0.0000110"
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.55%  '
$c = $ws.Range("D20")
$c.Value = "'6.73"
$c.Style = "Normal"
$ws.Range("E20").Value = '  +8.72%  '
$c = $ws.Range("D21")
$c.Value = "'72.22"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.36%  '
$c = $ws.Range("D22")
$c.Value = "'2.45"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.71%  '
$c = $ws.Range("D23")
$c.Value = "'9.84"
$c.Style = "Normal"
$ws.Range("E23").Value = '  +8.58%  '
$c = $ws.Range("D24")
$c.Value = "'232.39"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.32%  '
$c = $ws.Range("D25")
$c.Value = "'2.79"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +4.83%  '
$ws.Range("E26").Value = '  +0.03%  '
$c = $ws.Range("D27")
$c.Value = "'11.63"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.24%  '
$c = $ws.Range("D28")
$c.Value = "'41.72"
$c.Style = "Normal"
$ws.Range("E28").Value = '  +6.12%  '
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("E30").Value = '  -0.52%  '
$c = $ws.Range("D31")
$c.Value = "'175.36"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.19%  '
$c = $ws.Range("D32")
$c.Value = "'21.50"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.02%  '
$c = $ws.Range("D33")
$c.Value = "'0.0923"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.62%  '
$c = $ws.Range("D34")
$c.Value = "'5.65"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.66%  '
$c = $ws.Range("D35")
$c.Value = "'0.128"
$c.Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '
$c = $ws.Range("D36")
$c.Value = "'4.68"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +5.13%  '
$c = $ws.Range("D37")
$c.Value = "'0.0363"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +2.89%  '
$c = $ws.Range("D39")
$c.Value = "'3.86"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +13.71%  '
$ws.Range("E40").Value = '  +2.67%  '
$ws.Range("E41").Value = '  +0.85%  '
$c = $ws.Range("D42")
$c.Value = "'73.14"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +13.24%  '
$c = $ws.Range("D43")
$c.Value = "'13.65"
$c.Style = "Normal"
$ws.Range("E43").Value = '  +10.94%  '
$c = $ws.Range("D44")
$c.Value = "'6.33"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +22.63%  '
$ws.Range("E45").Value = '  +0.12%  '
$c = $ws.Range("D46")
$c.Value = "'1.38"
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.65%  '
$c = $ws.Range("D47")
$c.Value = "'8.75"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.46%  '
$c = $ws.Range("D48")
$c.Value = "'102.88"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +5.33%  '
$c = $ws.Range("D49")
$c.Value = "'0.0996"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("E50").Value = '  +2.82%  '
$c = $ws.Range("D51")
$c.Value = "'0.461"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +7.63%  '
